$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "1 Central Difference"
$ws.Range("B2").NumberFormat = "0%"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Report will be long"

# Row 3
$ws.Range("A3").Value = "2 Coordinates"
$ws.Range("B3").Value = 0.5

# Row 4
$ws.Range("A4").Value = "3 Central vs richardson"
$ws.Range("C4").Value = ""

# Row 5
$ws.Range("A5").Value = "4 Romberg"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = ""

# Row 7
$ws.Range("B7").NumberFormat = "0%"
$ws.Range("B7").Value = 1

# Column A width
$ws.Columns("A").ColumnWidth = 34.14

# Selection
$ws.Range("C7").Select() | Out-Null
